$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AV2").Value = 0
$ws.Range("AV2").Interior.Color = 255
$c = $ws.Range("AW2")
$c.Value = "'4152"
$c.Style = "Normal"

# Row 5
$c = $ws.Range("AW5")
$c.Value = "'2751"
$c.Style = "Normal"

# Row 18
$ws.Range("AV18").Value = 11
$ws.Range("AV18").Interior.Color = 65535
$c = $ws.Range("AW18")
$c.Value = "'3989"
$c.Style = "Normal"

# Row 20
$ws.Range("AV20").Value = 16
$ws.Range("AV20").Interior.Color = 65535
$c = $ws.Range("AW20")
$c.Value = "'4140"
$c.Style = "Normal"

# Row 21
$c = $ws.Range("AW21")
$c.Value = "'3137"
$c.Style = "Normal"

# Row 22
$c = $ws.Range("AW22")
$c.Value = "'4676"
$c.Style = "Normal"

# Row 23
$c = $ws.Range("AW23")
$c.Value = "'5144"
$c.Style = "Normal"

# Row 24
$ws.Range("AV24").Value = 33
$ws.Range("AV24").Interior.Color = 32768
$c = $ws.Range("AW24")
$c.Value = "'4601"
$c.Style = "Normal"

# Row 30
$ws.Range("AV30").Value = 20
$ws.Range("AV30").Interior.Color = 16777215
$c = $ws.Range("AW30")
$c.Value = "'4269"
$c.Style = "Normal"

# Row 31
$ws.Range("AV31").Value = 31
$ws.Range("AV31").Interior.Color = 32768
$c = $ws.Range("AW31")
$c.Value = "'4667"
$c.Style = "Normal"

# Row 32
$c = $ws.Range("AW32")
$c.Value = "'2624"
$c.Style = "Normal"

# Row 36
$c = $ws.Range("AW36")
$c.Value = "'3009"
$c.Style = "Normal"

# Row 38
$ws.Range("AV38").Value = 0
$ws.Range("AV38").Interior.Color = 255
$c = $ws.Range("AW38")
$c.Value = "'4676"
$c.Style = "Normal"

# Row 39
$c = $ws.Range("AW39")
$c.Value = "'4383"
$c.Style = "Normal"

# Row 41
$c = $ws.Range("AW41")
$c.Value = "'4169"
$c.Style = "Normal"

# Row 42
$c = $ws.Range("AW42")
$c.Value = "'2819"
$c.Style = "Normal"

# Row 46
$ws.Range("AV46").Value = 10
$ws.Range("AV46").Interior.Color = 65535
$c = $ws.Range("AW46")
$c.Value = "'4170"
$c.Style = "Normal"

# Row 47
$ws.Range("AV47").Value = 30
$ws.Range("AV47").Interior.Color = 16777215
$c = $ws.Range("AW47")
$c.Value = "'4994"
$c.Style = "Normal"

# Row 49
$ws.Range("AV49").Value = 27
$ws.Range("AV49").Interior.Color = 16777215
$c = $ws.Range("AW49")
$c.Value = "'4612"
$c.Style = "Normal"

# Row 50
$ws.Range("AV50").Value = 24
$ws.Range("AV50").Interior.Color = 16777215
$c = $ws.Range("AW50")
$c.Value = "'4712"
$c.Style = "Normal"

# Row 52
$c = $ws.Range("AW52")
$c.Value = "'4746"
$c.Style = "Normal"

# Row 53
$ws.Range("AV53").Value = 0
$ws.Range("AV53").Interior.Color = 255
$c = $ws.Range("AW53")
$c.Value = "'3336"
$c.Style = "Normal"

# Row 55
$c = $ws.Range("AW55")
$c.Value = "'3458"
$c.Style = "Normal"

# Row 56
$ws.Range("AV56").Value = 30
$ws.Range("AV56").Interior.Color = 16777215
$c = $ws.Range("AW56")
$c.Value = "'5008"
$c.Style = "Normal"

# Row 57
$ws.Range("AV57").Value = 20
$ws.Range("AV57").Interior.Color = 16777215
$c = $ws.Range("AW57")
$c.Value = "'4098"
$c.Style = "Normal"

# Row 58
$c = $ws.Range("AW58")
$c.Value = "'4085"
$c.Style = "Normal"

# Row 59
$ws.Range("AV59").Value = 0
$ws.Range("AV59").Interior.Color = 255
$c = $ws.Range("AW59")
$c.Value = "'3988"
$c.Style = "Normal"

# Row 60
$c = $ws.Range("AW60")
$c.Value = "'4193"
$c.Style = "Normal"

# Row 62
$c = $ws.Range("AW62")
$c.Value = "'3913"
$c.Style = "Normal"

# Row 63
$ws.Range("AV63").Value = 17
$ws.Range("AV63").Interior.Color = 65535
$c = $ws.Range("AW63")
$c.Value = "'3985"
$c.Style = "Normal"

# Row 64
$c = $ws.Range("AW64")
$c.Value = "'4181"
$c.Style = "Normal"

# Row 73
$c = $ws.Range("AW73")
$c.Value = "'2656"
$c.Style = "Normal"

# Row 76
$c = $ws.Range("AW76")
$c.Value = "'2636"
$c.Style = "Normal"

# Row 77
$c = $ws.Range("AW77")
$c.Value = "'2612"
$c.Style = "Normal"

# Row 96
$c = $ws.Range("AW96")
$c.Value = "'2470"
$c.Style = "Normal"

# Row 115
$c = $ws.Range("AW115")
$c.Value = "'5053"
$c.Style = "Normal"

# Row 117
$ws.Range("AV117").Value = 33
$ws.Range("AV117").Interior.Color = 32768
$c = $ws.Range("AW117")
$c.Value = "'5917"
$c.Style = "Normal"

# Row 118
$ws.Range("AV118").Value = 20
$ws.Range("AV118").Interior.Color = 16777215
$c = $ws.Range("AW118")
$c.Value = "'3233"
$c.Style = "Normal"

# Row 119
$c = $ws.Range("AW119")
$c.Value = "'1530"
$c.Style = "Normal"

# Row 125
$c = $ws.Range("AW125")
$c.Value = "'2148"
$c.Style = "Normal"

# Row 129
$c = $ws.Range("AW129")
$c.Value = "'2538"
$c.Style = "Normal"

# Row 132
$ws.Range("AV132").Value = 8
$ws.Range("AV132").Interior.Color = 65535
$c = $ws.Range("AW132")
$c.Value = "'4095"
$c.Style = "Normal"

# Row 133
$c = $ws.Range("AW133")
$c.Value = "'2481"
$c.Style = "Normal"

# Row 136
$ws.Range("AV136").Value = 39
$ws.Range("AV136").Interior.Color = 32768
$c = $ws.Range("AW136")
$c.Value = "'5865"
$c.Style = "Normal"

# Row 137
$ws.Range("B137").Value = '"L ᶻᵍˣ"'
$ws.Range("AV137").Value = 30
$ws.Range("AV137").Interior.Color = 16777215
$c = $ws.Range("AW137")
$c.Value = "'5409"
$c.Style = "Normal"

# Row 139
$ws.Range("AV139").Value = 34
$ws.Range("AV139").Interior.Color = 32768
$c = $ws.Range("AW139")
$c.Value = "'5845"
$c.Style = "Normal"

# Row 140
$c = $ws.Range("AW140")
$c.Value = "'2038"
$c.Style = "Normal"

# Row 142
$ws.Range("AV142").Value = 0
$ws.Range("AV142").Interior.Color = 255
$c = $ws.Range("AW142")
$c.Value = "'2924"
$c.Style = "Normal"

# Row 144
$ws.Range("AV144").Value = 6
$ws.Range("AV144").Interior.Color = 65535
$c = $ws.Range("AW144")
$c.Value = "'1621"
$c.Style = "Normal"

# Row 145
$c = $ws.Range("AW145")
$c.Value = "'1801"
$c.Style = "Normal"

# Row 147
$c = $ws.Range("AW147")
$c.Value = "'4483"
$c.Style = "Normal"

# Row 148
$c = $ws.Range("AW148")
$c.Value = "'2007"
$c.Style = "Normal"

# Row 150
$ws.Range("AV150").Value = 13
$ws.Range("AV150").Interior.Color = 65535
$c = $ws.Range("AW150")
$c.Value = "'3514"
$c.Style = "Normal"

# Row 151
$ws.Range("AV151").Value = 0
$ws.Range("AV151").Interior.Color = 255
$c = $ws.Range("AW151")
$c.Value = "'2296"
$c.Style = "Normal"

# Row 152
$ws.Range("AV152").Value = 7
$ws.Range("AV152").Interior.Color = 65535
$c = $ws.Range("AW152")
$c.Value = "'3975"
$c.Style = "Normal"

# Row 153
$c = $ws.Range("AW153")
$c.Value = "'1823"
$c.Style = "Normal"
